$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3, shifting the old row 3 down to row 4
$ws.Rows.Item(3).Insert()

# Update A2 with the new text
$ws.Range("A2").Value = "IB USer Demographic as at at -12-Dec-2024"

# The old A2 text now belongs in the newly inserted A3
$ws.Range("A3").Value = "IB User with blank Card via EBS (Tagged to Credit Card Brn) at at -12-Dec-2024"
